$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("A4").Value = "Tử Phủ Vũ Tướng"
$ws.Range("B4").Value = "Bạn là người có tính cách thích lãnh đạo, chỉ đạo."

# Row 5
$ws.Range("A5").Value = "Sát Phá Tham"
$ws.Range("B5").Value = "Bạn là người thiên về kinh doanh thương mại. "
$ws.Range("C5").Value = "Bản tính có xu hướng sát phạt, hơn thua, thích thay cũ đổi mới và có nhiều ham muốn."

# Row 6
$ws.Range("A6").Value = "Sát Phá Tham"

# Row 7
$ws.Range("A7").Value = "Tham Hỏa Linh"

# Row 8
$ws.Range("A8").Value = "Cơ Nguyệt Đồng Lương"
$ws.Range("B8").Value = "Bạn là tuýp người nhẹ nhàng, có nội tâm sâu sắc, phù hợp với môi trường giáo dục, công việc cần chuyên lý thuyết và tư duy cao."

# Update the selected cell to match the final saved view state
$ws.Range("G16").Select()
